# Apply the "Updated symbol list" price/volume refresh (Fri Dec 16 2022).
# Column D ("Price") cells are stored as literal text (t="inlineStr") in the
# workbook, so each target cell is forced to Text format before its value is
# written -- otherwise Excel would silently reinterpret a string such as
# "24.40" or "0.05910" as a number and drop the significant trailing zero.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.00"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.40"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.944"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05910"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.429"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.527"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.337"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7957"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07784"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03310"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03022"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09256"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.567"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001651"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04774"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006030"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006228"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005570"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001066"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001498"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.722"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3352"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1254"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006476"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04406"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007030"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1069"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003358"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009992"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002460"
$ws.Range("E45").Value = "44ACDXExchangeACXT"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005893"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9903"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1069"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
